$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 56168
$ws.Cells.Item(43, 10).Value = 96794
$ws.Cells.Item(43, 12).Value = 96794
$ws.Cells.Item(43, 14).Value = -96932
$ws.Cells.Item(64, 8).Value = 6999.75
$ws.Cells.Item(64, 9).Value = 6999.75
$ws.Cells.Item(64, 11).Value = 6999.75
$ws.Cells.Item(64, 13).Value = -6751.75
$ws.Cells.Item(67, 8).Value = 6999.75
$ws.Cells.Item(67, 9).Value = 6999.75
$ws.Cells.Item(67, 11).Value = 6999.75
$ws.Cells.Item(67, 13).Value = -6141.75
$ws.Cells.Item(99, 8).Value = 1011.25
$ws.Cells.Item(99, 9).Value = 816
$ws.Cells.Item(99, 10).Value = 1597
$ws.Cells.Item(99, 11).Value = 2448
$ws.Cells.Item(99, 12).Value = 4791
$ws.Cells.Item(99, 13).Value = -950
$ws.Cells.Item(99, 14).Value = -7787
$ws.Cells.Item(116, 8).Value = 5004.5
$ws.Cells.Item(116, 9).Value = 3698
$ws.Cells.Item(116, 11).Value = 3698
$ws.Cells.Item(116, 13).Value = -256
$ws.Cells.Item(132, 8).Value = 5297.857
$ws.Cells.Item(132, 9).Value = 5297.857
$ws.Cells.Item(132, 11).Value = 15893.571
$ws.Cells.Item(132, 13).Value = -13363.571
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5107.6665
$ws.Cells.Item(32, 9).Value = 5107.6665
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 5107.6665
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -4820.6665
$ws.Cells.Item(32, 14).Value = $null
$ws.Cells.Item(45, 8).Value = 1999.8
$ws.Cells.Item(45, 9).Value = 1999.8
$ws.Cells.Item(45, 11).Value = 1999.8
$ws.Cells.Item(45, 13).Value = -1622.8
$ws.Cells.Item(61, 8).Value = 10999.5
$ws.Cells.Item(61, 9).Value = 10999.5
$ws.Cells.Item(61, 11).Value = 10999.5
$ws.Cells.Item(61, 13).Value = -10787.5
$ws.Cells.Item(63, 8).Value = 2227.8333
$ws.Cells.Item(63, 9).Value = 1299
$ws.Cells.Item(63, 11).Value = 1299
$ws.Cells.Item(63, 13).Value = -613
$ws.Cells.Item(66, 8).Value = 2227.8333
$ws.Cells.Item(66, 9).Value = 1299
$ws.Cells.Item(66, 11).Value = 6495
$ws.Cells.Item(66, 13).Value = -3063
$ws.Cells.Item(88, 8).Value = 3316.2727
$ws.Cells.Item(88, 9).Value = 870
$ws.Cells.Item(88, 11).Value = 870
$ws.Cells.Item(88, 13).Value = -464
$ws.Cells.Item(91, 8).Value = 3316.2727
$ws.Cells.Item(91, 9).Value = 870
$ws.Cells.Item(91, 11).Value = 870
$ws.Cells.Item(91, 13).Value = 534
$ws.Cells.Item(136, 8).Value = 10999.5
$ws.Cells.Item(136, 9).Value = 10999.5
$ws.Cells.Item(136, 11).Value = 32998.5
$ws.Cells.Item(136, 13).Value = -30448.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 159.5
$ws.Cells.Item(7, 9).Value = 39.25
$ws.Cells.Item(7, 10).Value = 400
$ws.Cells.Item(7, 11).Value = 39.25
$ws.Cells.Item(7, 12).Value = 400
$ws.Cells.Item(7, 13).Value = 73.75
$ws.Cells.Item(7, 14).Value = -626
$ws.Cells.Item(15, 8).Value = 700
$ws.Cells.Item(15, 10).Value = 700
$ws.Cells.Item(15, 12).Value = 700
$ws.Cells.Item(15, 14).Value = -1154
$ws.Cells.Item(22, 8).Value = 10284.333
$ws.Cells.Item(22, 9).Value = 351
$ws.Cells.Item(22, 11).Value = 351
$ws.Cells.Item(22, 13).Value = -178
$ws.Cells.Item(33, 8).Value = 10000
$ws.Cells.Item(33, 9).Value = 10000
$ws.Cells.Item(33, 11).Value = 10000
$ws.Cells.Item(33, 13).Value = -9664
$ws.Cells.Item(94, 8).Value = 2543.8667
$ws.Cells.Item(94, 9).Value = 2543.8667
$ws.Cells.Item(94, 11).Value = 2543.8667
$ws.Cells.Item(94, 13).Value = -2092.8667
$ws.Cells.Item(134, 8).Value = 9399.333000000001
$ws.Cells.Item(134, 9).Value = 9719.299999999999
$ws.Cells.Item(134, 11).Value = 29157.9
$ws.Cells.Item(134, 13).Value = -26622.9
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 19266.5
$ws.Cells.Item(28, 10).Value = 19266.5
$ws.Cells.Item(28, 12).Value = 19266.5
$ws.Cells.Item(28, 14).Value = -19756.5
$ws.Cells.Item(31, 8).Value = 2362.4
$ws.Cells.Item(31, 9).Value = 2236
$ws.Cells.Item(31, 11).Value = 2236
$ws.Cells.Item(31, 13).Value = -1941
$ws.Cells.Item(34, 8).Value = 2362.4
$ws.Cells.Item(34, 9).Value = 2236
$ws.Cells.Item(34, 11).Value = 2236
$ws.Cells.Item(34, 13).Value = -2034
$ws.Cells.Item(62, 8).Value = 2998
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 14).Value = $null
$ws.Cells.Item(65, 8).Value = 2998
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 14).Value = $null
$ws.Cells.Item(74, 8).Value = 40000
$ws.Cells.Item(74, 10).Value = 40000
$ws.Cells.Item(74, 12).Value = 40000
$ws.Cells.Item(74, 14).Value = -41748
$ws.Cells.Item(77, 8).Value = 40000
$ws.Cells.Item(77, 10).Value = 40000
$ws.Cells.Item(77, 12).Value = 120000
$ws.Cells.Item(77, 14).Value = -128736
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 13).Value = $null
$ws.Cells.Item(134, 8).Value = 7133.6665
$ws.Cells.Item(134, 9).Value = 6260.5
$ws.Cells.Item(134, 11).Value = 18781.5
$ws.Cells.Item(134, 13).Value = -16246.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1483.5588
$ws.Cells.Item(4, 10).Value = 1518
$ws.Cells.Item(4, 12).Value = 4554
$ws.Cells.Item(4, 14).Value = -4778
$ws.Cells.Item(107, 8).Value = 579.73334
$ws.Cells.Item(107, 9).Value = 483.2
$ws.Cells.Item(107, 10).Value = 772.8
$ws.Cells.Item(107, 11).Value = 1449.6
$ws.Cells.Item(107, 12).Value = 2318.4
$ws.Cells.Item(107, 13).Value = 470.4000000000001
$ws.Cells.Item(107, 14).Value = -6158.4
$ws.Cells.Item(113, 8).Value = 543.2353000000001
$ws.Cells.Item(113, 10).Value = 673.125
$ws.Cells.Item(113, 12).Value = 2019.375
$ws.Cells.Item(113, 14).Value = -6359.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(100, 8).Value = 22500
$ws.Cells.Item(100, 10).Value = 22500
$ws.Cells.Item(100, 12).Value = 22500
$ws.Cells.Item(100, 14).Value = -24664
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2361.6365
$ws.Cells.Item(7, 9).Value = 2131.5557
$ws.Cells.Item(7, 11).Value = 2131.5557
$ws.Cells.Item(7, 13).Value = -2019.5557
$ws.Cells.Item(16, 8).Value = 2937.3333
$ws.Cells.Item(16, 9).Value = 1514.1428
$ws.Cells.Item(16, 11).Value = 1514.1428
$ws.Cells.Item(16, 13).Value = -1344.1428
$ws.Cells.Item(93, 8).Value = 968
$ws.Cells.Item(93, 9).Value = 900
$ws.Cells.Item(93, 11).Value = 900
$ws.Cells.Item(93, 13).Value = 348
$ws.Cells.Item(126, 8).Value = 2361.6365
$ws.Cells.Item(126, 9).Value = 2131.5557
$ws.Cells.Item(126, 11).Value = 6394.6671
$ws.Cells.Item(126, 13).Value = -3924.6671
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 14).Value = $null
$ws.Cells.Item(100, 8).Value = 848.5
$ws.Cells.Item(100, 10).Value = 898
$ws.Cells.Item(100, 12).Value = 1796
$ws.Cells.Item(100, 14).Value = -2878
$ws.Cells.Item(126, 8).Value = 2128.4707
$ws.Cells.Item(126, 10).Value = 3123.75
$ws.Cells.Item(126, 12).Value = 9371.25
$ws.Cells.Item(126, 14).Value = -14311.25
$ws.Cells.Item(132, 8).Value = 1104.5
$ws.Cells.Item(132, 9).Value = 1104.5
$ws.Cells.Item(132, 11).Value = 3313.5
$ws.Cells.Item(132, 13).Value = -783.5
